$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new column at N (ChlA/TP/DOC block gets a new "SW DOC" column) ---
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Header for the newly inserted column
$ws.Range("N2").Value = "SW DOC (g/m3)"

# --- Updated Secchi values (column K): new 2-decimal number format, no border, ---
# --- vertically centered / default horizontal alignment (applied to K3:K7 as one block) ---
$ws.Range("K3").Value = 4.4680929999999996
$ws.Range("K4").Value = 3.0271870000000001
$ws.Range("K5").Value = 4.5320900000000002
$ws.Range("K6").Value = 5.3249639999999996
$ws.Range("K7").Value = 4.4680929999999996

$ws.Range("K3:K7").NumberFormat = "0.00"
$ws.Range("K3:K7").HorizontalAlignment = 1
$ws.Range("K3:K7").VerticalAlignment = -4108
$ws.Range("K7").Borders.LineStyle = -4142

# --- Updated ChlA (L), TP (M) and new SW DOC (N) values, 2-decimal number format ---
# (keeps each row's existing alignment / border, matches the per-row formatting already present)
$ws.Range("L3").Value = 2.0666009999999999
$ws.Range("M3").Value = 6.2953239999999999
$ws.Range("N3").Value = 9.3538669999999993

$ws.Range("L4").Value = 9.2059099999999994
$ws.Range("M4").Value = 135.56549999999999
$ws.Range("N4").Value = 5.0944419999999999

$ws.Range("L5").Value = 1.4116029999999999
$ws.Range("M5").Value = 3.787623
$ws.Range("N5").Value = 7.7224019999999998

$ws.Range("L6").Value = 2.234003
$ws.Range("M6").Value = 10.98129
$ws.Range("N6").Value = 5.1144499999999997

$ws.Range("L7").Value = 2.0666009999999999
$ws.Range("M7").Value = 6.2953239999999999
$ws.Range("N7").Value = 9.3538669999999993

$ws.Range("L3:N7").NumberFormat = "0.00"

# --- Match the workbook's recorded selection after the edit ---
[void]$ws.Range("O8").Select()
